## Add new rows 54-69 to the symbol table on Sheet1, matching the new
## sharedStrings entries, and tweak the sheet view (scroll position/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New symbol-table rows --------------------------------------------------

$ws.Range("A54").Value = "`$SSB_{MSY}^P`$"
$ws.Range("B54").Value = "predator spawning stock biomass at maximum sustainable yield"

$ws.Range("A55").Value = "`$GR`$"
$ws.Range("B55").Value = "gross revenue"

$ws.Range("A56").Value = "`$NR`$"
$ws.Range("B56").Value = "net operating revenues"

$ws.Range("A57").Value = "`$i`$"
# B57 uses rich text: "fleet (trawl, " + italic "t, " + "or purse seine, " + italic "s" + ")"
$ws.Range("B57").Value = "fleet (trawl, t, or purse seine, s)"
$b57 = $ws.Range("B57")
$b57.Characters(15, 3).Font.Italic = $true
$b57.Characters(18, 16).Font.Italic = $false
$b57.Characters(34, 1).Font.Italic = $true
$b57.Characters(35, 1).Font.Italic = $false

$ws.Range("A58").Value = "`$q`$"
$ws.Range("B58").Value = "quantity landed"

$ws.Range("A59").Value = "`$c`$"
$ws.Range("B59").Value = "cost function"

$ws.Range("A60").Value = "`$p`$"
$ws.Range("B60").Value = "function relating landings to prices"

$ws.Range("A61").Value = "`$t`$"
$ws.Range("B61").Value = "denotes trawl fishery in economic model"

$ws.Range("A62").Value = "`$s`$"
$ws.Range("B62").Value = "denotes purse seine fishery in economic model"

$ws.Range("A63").Value = "`$a_i`$"
$ws.Range("B63").Value = "parameter of price and landing model"

$ws.Range("A64").Value = "`$b`$"
$ws.Range("B64").Value = "parameter of price and landing model"

$ws.Range("A65").Value = "`$\gamma`$"
$ws.Range("B65").Value = "parameter of price and landing model"

$ws.Range("A66").Value = "`$\alpha`$"
$ws.Range("B66").Value = "parameter of price and landing model"

$ws.Range("A67").Value = "`$\theta_i`$"
$ws.Range("B67").Value = "parameter of price and landing model"

$ws.Range("A68").Value = "`$\beta`$"
$ws.Range("B68").Value = "parameter of price and landing model"

$ws.Range("A69").Value = "`$\xi`$"
$ws.Range("B69").Value = "parameter of economic stationarity metric"

# --- View tweaks -------------------------------------------------------------

$ws.Range("E43").Select()
